# Fruta / hortaliza, semanal
#
# A new weekly price-report record for "Repollo" (Crespo record, Primera)
# at the Macroferia Regional de Talca is inserted as row 470. Excel's
# native row-insert shifts the existing rows 470-529 down to 471-530
# (dimension grows from A1:R529 to A1:R530), and the new row's values are
# then populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 470, pushing the old rows 470:529 down to 471:530.
$ws.Rows("470:470").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A470").Value = 5
$ws.Range("B470").Value = "Macroferia Regional de Talca"
$ws.Range("C470").Value = "Maule"
$ws.Range("D470").Value = 45154
$ws.Range("E470").Value = 7
$ws.Range("F470").Value = 100112006
$ws.Range("G470").Value = "Repollo"
$ws.Range("H470").Value = "Crespo record"
$ws.Range("I470").Value = "Primera"
$ws.Range("J470").Value = 5000
$ws.Range("K470").Value = 600
$ws.Range("L470").Value = 600
$ws.Range("M470").Value = 600
$ws.Range("N470").Value = "$/unidad"
$ws.Range("O470").Value = "Región del Maule"
$ws.Range("P470").Value = 600
$ws.Range("Q470").Value = 1
$ws.Range("R470").Value = "Hortaliza"
